# Auto-generated Excel COM-interop edit script.
# Updates specific LeveProfit / average-price cells across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets of the Sargatanas Profits workbook, reflecting a refresh
# of market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 124
$ws.Range("H124").Value = 102000
$ws.Range("J124").Value = 102000
$ws.Range("L124").Value = 102000
$ws.Range("N124").Value = -111820

# Row 127
$ws.Range("H127").Value = 4437.5
$ws.Range("I127").Value = 4437.5
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 13312.5
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -8352.5
$ws.Range("N127").Value = ""

# Row 132
$ws.Range("H132").Value = 2329.1692
$ws.Range("I132").Value = 2207.2166
$ws.Range("J132").Value = 3792.6
$ws.Range("K132").Value = 6621.649800000001
$ws.Range("L132").Value = 11377.8
$ws.Range("M132").Value = -4091.649800000001
$ws.Range("N132").Value = -16437.8

# Row 137
$ws.Range("H137").Value = 4342.9287
$ws.Range("I137").Value = 20000
$ws.Range("K137").Value = 60000
$ws.Range("M137").Value = -57450

# Row 138
$ws.Range("H138").Value = 6200.095
$ws.Range("I138").Value = 1674.6666
$ws.Range("J138").Value = 8010.2666
$ws.Range("K138").Value = 5023.9998
$ws.Range("L138").Value = 24030.7998
$ws.Range("M138").Value = 116.0002000000004
$ws.Range("N138").Value = -34310.7998

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = ""

# Row 40
$ws.Range("H40").Value = 2231401.2
$ws.Range("I40").Value = 11444.182
$ws.Range("J40").Value = 8336283
$ws.Range("K40").Value = 11444.182
$ws.Range("L40").Value = 8336283
$ws.Range("M40").Value = -11269.182
$ws.Range("N40").Value = -8336633

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 4710.8335
$ws.Range("I102").Value = 4567.625
$ws.Range("K102").Value = 4567.625
$ws.Range("M102").Value = -2945.625

# Row 122
$ws.Range("H122").Value = 15608.5
$ws.Range("I122").Value = 18811.416
$ws.Range("K122").Value = 56434.24800000001
$ws.Range("M122").Value = -53984.24800000001

# Row 132
$ws.Range("H132").Value = 6914.6562
$ws.Range("I132").Value = 3471.6667
$ws.Range("K132").Value = 10415.0001
$ws.Range("M132").Value = -7885.000100000001

# Row 136
$ws.Range("H136").Value = 11785.866
$ws.Range("I136").Value = 2798.8
$ws.Range("K136").Value = 8396.400000000001
$ws.Range("M136").Value = -5846.400000000001

# Row 32
$ws.Range("H32").Value = 1840650.4
$ws.Range("I32").Value = 1986372.8
$ws.Range("K32").Value = 1986372.8
$ws.Range("M32").Value = -1986085.8

# Row 57
$ws.Range("H57").Value = 5500
$ws.Range("I57").Value = 5500
$ws.Range("K57").Value = 5500
$ws.Range("M57").Value = -5016

# Row 61
$ws.Range("H61").Value = 11785.866
$ws.Range("I61").Value = 2798.8
$ws.Range("K61").Value = 2798.8
$ws.Range("M61").Value = -2586.8

# Row 74
$ws.Range("H74").Value = 66509.16
$ws.Range("I74").Value = 108099.53
$ws.Range("J74").Value = 4123.6
$ws.Range("K74").Value = 108099.53
$ws.Range("L74").Value = 4123.6
$ws.Range("M74").Value = -107225.53
$ws.Range("N74").Value = -5871.6

# Row 77
$ws.Range("H77").Value = 66509.16
$ws.Range("I77").Value = 108099.53
$ws.Range("J77").Value = 4123.6
$ws.Range("K77").Value = 540497.65
$ws.Range("L77").Value = 20618
$ws.Range("M77").Value = -536129.65
$ws.Range("N77").Value = -29354

# Row 88
$ws.Range("H88").Value = 1818.2632
$ws.Range("J88").Value = 2002.6
$ws.Range("L88").Value = 2002.6
$ws.Range("N88").Value = -2814.6

# Row 91
$ws.Range("H91").Value = 1818.2632
$ws.Range("J91").Value = 2002.6
$ws.Range("L91").Value = 2002.6
$ws.Range("N91").Value = -4810.6

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2955.88
$ws.Range("I105").Value = 2212
$ws.Range("J105").Value = 4278.3335
$ws.Range("K105").Value = 2212
$ws.Range("L105").Value = 4278.3335
$ws.Range("M105").Value = -465
$ws.Range("N105").Value = -7772.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 113
$ws.Range("H113").Value = 4889.727
$ws.Range("I113").Value = 3197.6
$ws.Range("J113").Value = 6299.8335
$ws.Range("K113").Value = 3197.6
$ws.Range("L113").Value = 6299.8335
$ws.Range("M113").Value = -1027.6
$ws.Range("N113").Value = -10639.8335

# Row 126
$ws.Range("H126").Value = 11217.667
$ws.Range("J126").Value = 8101.5
$ws.Range("L126").Value = 24304.5
$ws.Range("N126").Value = -29244.5

# Row 134
$ws.Range("H134").Value = 8603.219999999999
$ws.Range("I134").Value = 10170.25
$ws.Range("J134").Value = 7600.32
$ws.Range("K134").Value = 30510.75
$ws.Range("L134").Value = 22800.96
$ws.Range("M134").Value = -27975.75
$ws.Range("N134").Value = -27870.96

# Row 16
$ws.Range("H16").Value = 4889.727
$ws.Range("I16").Value = 3197.6
$ws.Range("J16").Value = 6299.8335
$ws.Range("K16").Value = 3197.6
$ws.Range("L16").Value = 6299.8335
$ws.Range("M16").Value = -2910.6
$ws.Range("N16").Value = -6873.8335

# Row 31
$ws.Range("H31").Value = 6667.6064
$ws.Range("I31").Value = 2954.4412
$ws.Range("J31").Value = 11343.444
$ws.Range("K31").Value = 2954.4412
$ws.Range("L31").Value = 11343.444
$ws.Range("M31").Value = -2659.4412
$ws.Range("N31").Value = -11933.444

# Row 34
$ws.Range("H34").Value = 6667.6064
$ws.Range("I34").Value = 2954.4412
$ws.Range("J34").Value = 11343.444
$ws.Range("K34").Value = 2954.4412
$ws.Range("L34").Value = 11343.444
$ws.Range("M34").Value = -2752.4412
$ws.Range("N34").Value = -11747.444

# Row 99
$ws.Range("H99").Value = 11217.667
$ws.Range("J99").Value = 8101.5
$ws.Range("L99").Value = 8101.5
$ws.Range("N99").Value = -11097.5

$ws = $wb.Worksheets.Item("CUL")
# Row 111
$ws.Range("H111").Value = 15000.8
$ws.Range("I111").Value = 15000.8
$ws.Range("K111").Value = 45002.39999999999
$ws.Range("M111").Value = -41935.39999999999

# Row 113
$ws.Range("H113").Value = 2789.0952
$ws.Range("I113").Value = 1604.25
$ws.Range("K113").Value = 4812.75
$ws.Range("M113").Value = -2642.75

# Row 119
$ws.Range("H119").Value = 4698.5
$ws.Range("I119").Value = 4698.5
$ws.Range("K119").Value = 14095.5
$ws.Range("M119").Value = -9257.5

# Row 122
$ws.Range("H122").Value = 3145458
$ws.Range("I122").Value = 5658705.5
$ws.Range("J122").Value = 3898.5
$ws.Range("K122").Value = 50928349.5
$ws.Range("L122").Value = 35086.5
$ws.Range("M122").Value = -50925899.5
$ws.Range("N122").Value = -39986.5

# Row 98
$ws.Range("H98").Value = 461.6
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 3000
$ws.Range("N98").Value = -5996

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 6205.4
$ws.Range("I132").Value = 2756.75
$ws.Range("K132").Value = 8270.25
$ws.Range("M132").Value = -5740.25

# Row 39
$ws.Range("H39").Value = 50735.668
$ws.Range("J39").Value = 50735.668
$ws.Range("L39").Value = 50735.668
$ws.Range("N39").Value = -51799.668

# Row 80
$ws.Range("H80").Value = 4381.5713
$ws.Range("J80").Value = 4900.25
$ws.Range("L80").Value = 4900.25
$ws.Range("N80").Value = -6896.25

# Row 83
$ws.Range("H83").Value = 4381.5713
$ws.Range("J83").Value = 4900.25
$ws.Range("L83").Value = 24501.25
$ws.Range("N83").Value = -34485.25

# Row 97
$ws.Range("H97").Value = 1632.0541
$ws.Range("I97").Value = 1446.6957
$ws.Range("K97").Value = 1446.6957
$ws.Range("M97").Value = -950.6957

$ws = $wb.Worksheets.Item("LTW")
# Row 113
$ws.Range("H113").Value = 6369.4614
$ws.Range("I113").Value = 4667.6665
$ws.Range("K113").Value = 4667.6665
$ws.Range("M113").Value = -2497.6665

# Row 132
$ws.Range("H132").Value = 20842390
$ws.Range("I132").Value = 41671830
$ws.Range("K132").Value = 125015490
$ws.Range("M132").Value = -125012960

# Row 136
$ws.Range("H136").Value = 8799.909
$ws.Range("I136").Value = 4149.5
$ws.Range("J136").Value = 13176.765
$ws.Range("K136").Value = 12448.5
$ws.Range("L136").Value = 39530.295
$ws.Range("M136").Value = -9898.5
$ws.Range("N136").Value = -44630.295

# Row 22
$ws.Range("H22").Value = 26924.5
$ws.Range("I22").Value = 7724
$ws.Range("K22").Value = 7724
$ws.Range("M22").Value = -7429

# Row 27
$ws.Range("H27").Value = 26924.5
$ws.Range("I27").Value = 7724
$ws.Range("K27").Value = 7724
$ws.Range("M27").Value = -7617

# Row 46
$ws.Range("H46").Value = 20801340
$ws.Range("I46").Value = 34482760
$ws.Range("J46").Value = 18521102
$ws.Range("K46").Value = 34482760
$ws.Range("L46").Value = 18521102
$ws.Range("M46").Value = -34482572
$ws.Range("N46").Value = -18521478

# Row 61
$ws.Range("H61").Value = 6369.4614
$ws.Range("I61").Value = 4667.6665
$ws.Range("K61").Value = 4667.6665
$ws.Range("M61").Value = -4465.6665

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5844.95
$ws.Range("I122").Value = 4275.125
$ws.Range("K122").Value = 12825.375
$ws.Range("M122").Value = -10375.375

# Row 132
$ws.Range("H132").Value = 13532200
$ws.Range("I132").Value = 18523730
$ws.Range("K132").Value = 55571190
$ws.Range("M132").Value = -55568660

# Row 136
$ws.Range("H136").Value = 55560176
$ws.Range("I136").Value = 111112776
$ws.Range("K136").Value = 333338328
$ws.Range("M136").Value = -333335778

# Row 62
$ws.Range("H62").Value = 11492.714
$ws.Range("I62").Value = 11492.714
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 11492.714
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -10868.714
$ws.Range("N62").Value = ""

# Row 65
$ws.Range("H65").Value = 11492.714
$ws.Range("I65").Value = 11492.714
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 57463.57
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -54343.57
$ws.Range("N65").Value = ""

